$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 19

$ws.Range("A$row").Value = 18
$ws.Range("B$row").Value = "Saturday, Jan 14"
$ws.Range("C$row").Value = "8:00 PM"
$ws.Range("D$row").Value = "FR9890"
$ws.Range("E$row").Value = "Milan"
$ws.Range("F$row").Value = "(BGY)"
$ws.Range("G$row").Value = "Ryanair "
$ws.Range("H$row").Value = "B738"
$ws.Range("I$row").Value = "(9H-QEC)"
$ws.Range("J$row").Value = "8:11 PM"
$ws.Range("K$row").Font.Bold = $false
$ws.Range("L$row").Value = "0 hours, 11 minutes"
$ws.Range("M$row").Font.Bold = $false
